# Remove numero compte bancaire unique validation
# Collapse the data to a single record (NACER YASSINE / 389-AOURIR contract,
# trimestrielle periodicity, 8000/800/7200 amounts) plus its totals row,
# and delete the now-unused rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 8 (old extra data rows + old totals row) first,
# working from the bottom up so row indices stay stable.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Update remaining row 2 with the new single data record.
$ws.Range("A2").Value = "NACER YASSINE"
$ws.Range("B2").Value = "L234567"
# Leading apostrophe forces this long, purely-numeric-looking value to be
# stored as text instead of being coerced into a (precision-losing) number.
$ws.Range("C2").Value = "'78017098772736274634834384"
$ws.Range("D2").Value = "TOUHAMI"
$ws.Range("E2").Value = "ATTIJARI WAFA BANK"
$ws.Range("F2").Value = "Point de vente"
$ws.Range("G2").Value = "389/AOURIR"
$ws.Range("H2").Value = "trimestrielle"
$ws.Range("I2").Value = 8000
$ws.Range("J2").Value = 800
$ws.Range("K2").Value = 7200

# Update row 3 to become the totals row (blank labels, summed amounts).
$ws.Range("A3").Value = " "
$ws.Range("B3").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("I3").Value = 8000
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 7200
